$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" holds the note in A1 with the conversion rates text
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.55 = 17913.64 pesos`n✅ 17913.64 pesos = 4.52 = 919.52 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas" holds the numeric rate table updates
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 220
$wsTasas.Range("O10").Value = 3941
$wsTasas.Range("N12").Value = 3964.48
$wsTasas.Range("O12").Value = 203.5
